# "agregue cositas en la caja diaria y corregi un par de errores"
#
# - Adds two new columns to the daily-cash sheet: tipoEdad (H) and
#   tipoMordida (I).
# - Adds two new product rows (6 and 7) to the table.
#
# All the "data" columns (name/tipo/peso/price/stock/fecha, and the two
# new ones) are kept as plain text, even the ones that look like numbers
# or dates, exactly like the rest of the sheet already does -- only the
# "id" column (G) is a real number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new header cells --------------------------------------------------
$ws.Range("H1").Value = "tipoEdad"
$ws.Range("I1").Value = "tipoMordida"

# --- force the new text columns to keep their literal text -------------
# Without this, values such as "8", "10000" or "2023-09-07" would be
# auto-converted to numbers/dates by Excel's normal typing behaviour.
$textCells = $ws.Range("A6:F7,H7:I7")
$textCells.NumberFormat = "@"

# --- row 6 ---------------------------------------------------------------
$ws.Range("A6").Value = "NUTRIBON"
$ws.Range("B6").Value = "GATO"
$ws.Range("C6").Value = "8"
$ws.Range("D6").Value = "10000"
$ws.Range("E6").Value = "1"
$ws.Range("F6").Value = "2023-09-07"
$ws.Range("G6").Value = 5

# --- row 7 ---------------------------------------------------------------
$ws.Range("A7").Value = "NUTRIBON"
$ws.Range("B7").Value = "PERRO"
$ws.Range("C7").Value = "12"
$ws.Range("D7").Value = "123123"
$ws.Range("E7").Value = "22"
$ws.Range("F7").Value = "2023-09-07"
$ws.Range("G7").Value = 6
$ws.Range("H7").Value = "ADULTO"
$ws.Range("I7").Value = "GRANDE"

# --- drop the temporary "Text" number format again ---------------------
# so the new cells fall back onto the workbook's default (unstyled) look,
# matching every other text cell already on the sheet.
$textCells.Style = "Normal"
